$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (Hora) changes uniformly from 19 to 20 for all data rows
$gRange = $ws.Range("G2:G51")
$gRange.NumberFormat = "@"
$gRange.Value = "20"
$gRange.Style = "Normal"

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "285.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.30%"
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "5.14%"
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.955"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.43%"
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06694"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.87%"
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.335"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.21%"
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.364"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.74%"
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9412"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "5.21%"
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1579"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3.41%"
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06579"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "13.39%"
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07583"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.00%"
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02948"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.99%"
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09000"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.04%"
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001588"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.08%"
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04486"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.90%"
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006455"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.53%"
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006387"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "4.75%"
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.491"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.49%"
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.384"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.37%"
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.247"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.95%"
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3204"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.83%"
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1309"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-3.08%"
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.080"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "4.35%"
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1524"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.43%"
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001179"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.25%"
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004492"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.96%"
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "5.78%"
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001615"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-2.29%"
$ws.Range("E28").Style = "Normal"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04205"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.39%"
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006716"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.10%"
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1254"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-10.09%"
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.31%"
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01221"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "11.19%"
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005741"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.45%"
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "25.93%"
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01305"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-29.43%"
$ws.Range("E47").Style = "Normal"
